$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($row, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Insert two new rows before row 723, shifting the existing data
# (old rows 723-738) down to rows 725-740.
$ws.Rows.Item(723).Insert()
$ws.Rows.Item(723).Insert()

# Populate the two newly inserted rows with the new weekly data.
Set-RowValues 723 @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45239, 13, 100112013, 'Alcachofa', 'Española', 'Primera', 7000, 450, 480, 465, '$/unidad', 'Región Metropolitana', 465, 1, 'Hortaliza')
Set-RowValues 724 @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45239, 13, 100112013, 'Alcachofa', 'Española', 'Segunda', 5200, 320, 350, 335, '$/unidad', 'Región Metropolitana', 335, 1, 'Hortaliza')
